$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").Value = "'258.15"
$ws.Range("D2").Style = $origStyle
$origStyle = $ws.Range("E2").Style
$ws.Range("E2").Value = "'5.34%"
$ws.Range("E2").Style = $origStyle
$origStyle = $ws.Range("D3").Style
$ws.Range("D3").Value = "'27.81"
$ws.Range("D3").Style = $origStyle
$origStyle = $ws.Range("E3").Style
$ws.Range("E3").Value = "'-2.45%"
$ws.Range("E3").Style = $origStyle
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").Value = "'5.223"
$ws.Range("D4").Style = $origStyle
$origStyle = $ws.Range("E4").Style
$ws.Range("E4").Value = "'-0.13%"
$ws.Range("E4").Style = $origStyle
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").Value = "'0.05947"
$ws.Range("D5").Style = $origStyle
$origStyle = $ws.Range("E5").Style
$ws.Range("E5").Value = "'4.41%"
$ws.Range("E5").Style = $origStyle
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").Value = "'6.701"
$ws.Range("D6").Style = $origStyle
$origStyle = $ws.Range("E6").Style
$ws.Range("E6").Value = "'1.37%"
$ws.Range("E6").Style = $origStyle
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.8717"
$ws.Range("D7").Style = $origStyle
$origStyle = $ws.Range("E7").Style
$ws.Range("E7").Value = "'2.52%"
$ws.Range("E7").Style = $origStyle
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").Value = "'1.039"
$ws.Range("D8").Style = $origStyle
$origStyle = $ws.Range("E8").Style
$ws.Range("E8").Value = "'21.22%"
$ws.Range("E8").Style = $origStyle
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.1428"
$ws.Range("D9").Style = $origStyle
$origStyle = $ws.Range("E9").Style
$ws.Range("E9").Value = "'4.31%"
$ws.Range("E9").Style = $origStyle
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").Value = "'0.07241"
$ws.Range("D10").Style = $origStyle
$origStyle = $ws.Range("E10").Style
$ws.Range("E10").Value = "'2.92%"
$ws.Range("E10").Style = $origStyle
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.03225"
$ws.Range("D11").Style = $origStyle
$origStyle = $ws.Range("E11").Style
$ws.Range("E11").Value = "'2.82%"
$ws.Range("E11").Style = $origStyle
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.09234"
$ws.Range("D12").Style = $origStyle
$origStyle = $ws.Range("E12").Style
$ws.Range("E12").Value = "'0.15%"
$ws.Range("E12").Style = $origStyle
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").Value = "'0.001563"
$ws.Range("D13").Style = $origStyle
$origStyle = $ws.Range("E13").Style
$ws.Range("E13").Value = "'1.64%"
$ws.Range("E13").Style = $origStyle
$ws.Range("B14").Value = 'One'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").Value = "'0.0006070"
$ws.Range("D14").Style = $origStyle
$origStyle = $ws.Range("E14").Style
$ws.Range("E14").Value = "'-93.94%"
$ws.Range("E14").Style = $origStyle
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").Value = "'0.005970"
$ws.Range("D15").Style = $origStyle
$origStyle = $ws.Range("E15").Style
$ws.Range("E15").Value = "'0.87%"
$ws.Range("E15").Style = $origStyle
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").Value = "'3.486"
$ws.Range("D16").Style = $origStyle
$origStyle = $ws.Range("E16").Style
$ws.Range("E16").Value = "'-0.06%"
$ws.Range("E16").Style = $origStyle
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").Value = "'3.267"
$ws.Range("D17").Style = $origStyle
$origStyle = $ws.Range("E17").Style
$ws.Range("E17").Value = "'2.25%"
$ws.Range("E17").Style = $origStyle
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").Value = "'2.210"
$ws.Range("D18").Style = $origStyle
$origStyle = $ws.Range("E18").Style
$ws.Range("E18").Value = "'1.66%"
$ws.Range("E18").Style = $origStyle
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").Value = "'0.3148"
$ws.Range("D19").Style = $origStyle
$origStyle = $ws.Range("E19").Style
$ws.Range("E19").Value = "'-0.67%"
$ws.Range("E19").Style = $origStyle
$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").Value = "'0.03623"
$ws.Range("D20").Style = $origStyle
$origStyle = $ws.Range("E20").Style
$ws.Range("E20").Value = "'10.72%"
$ws.Range("E20").Style = $origStyle
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").Value = "'0.1291"
$ws.Range("D21").Style = $origStyle
$origStyle = $ws.Range("E21").Style
$ws.Range("E21").Value = "'0.29%"
$ws.Range("E21").Style = $origStyle
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").Value = "'3.532"
$ws.Range("D22").Style = $origStyle
$origStyle = $ws.Range("E22").Style
$ws.Range("E22").Value = "'1.34%"
$ws.Range("E22").Style = $origStyle
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").Value = "'0.04186"
$ws.Range("D23").Style = $origStyle
$origStyle = $ws.Range("E23").Style
$ws.Range("E23").Value = "'2.40%"
$ws.Range("E23").Style = $origStyle
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").Value = "'0.1399"
$ws.Range("D24").Style = $origStyle
$origStyle = $ws.Range("E24").Style
$ws.Range("E24").Value = "'1.44%"
$ws.Range("E24").Style = $origStyle
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").Value = "'0.001221"
$ws.Range("D25").Style = $origStyle
$origStyle = $ws.Range("E25").Style
$ws.Range("E25").Value = "'-0.02%"
$ws.Range("E25").Style = $origStyle
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").Value = "'0.004554"
$ws.Range("D26").Style = $origStyle
$origStyle = $ws.Range("E26").Style
$ws.Range("E26").Value = "'10.11%"
$ws.Range("E26").Style = $origStyle
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").Value = "'0.0001202"
$ws.Range("D27").Style = $origStyle
$origStyle = $ws.Range("E27").Style
$ws.Range("E27").Value = "'0.18%"
$ws.Range("E27").Style = $origStyle
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").Value = "'0.0001938"
$ws.Range("D28").Style = $origStyle
$origStyle = $ws.Range("E28").Style
$ws.Range("E28").Value = "'33.73%"
$ws.Range("E28").Style = $origStyle
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").Value = "'0.03822"
$ws.Range("D40").Style = $origStyle
$origStyle = $ws.Range("E40").Style
$ws.Range("E40").Value = "'1.54%"
$ws.Range("E40").Style = $origStyle
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").Value = "'0.005580"
$ws.Range("D41").Style = $origStyle
$origStyle = $ws.Range("E41").Style
$ws.Range("E41").Value = "'7.53%"
$ws.Range("E41").Style = $origStyle
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").Value = "'0.1108"
$ws.Range("D42").Style = $origStyle
$origStyle = $ws.Range("E42").Style
$ws.Range("E42").Value = "'4.20%"
$ws.Range("E42").Style = $origStyle
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").Value = "'0.002385"
$ws.Range("D43").Style = $origStyle
$origStyle = $ws.Range("E43").Style
$ws.Range("E43").Value = "'3.72%"
$ws.Range("E43").Style = $origStyle
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").Value = "'0.009949"
$ws.Range("D44").Style = $origStyle
$origStyle = $ws.Range("E44").Style
$ws.Range("E44").Value = "'8.89%"
$ws.Range("E44").Style = $origStyle
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").Value = "'0.00005430"
$ws.Range("D45").Style = $origStyle
$origStyle = $ws.Range("E45").Style
$ws.Range("E45").Value = "'2.81%"
$ws.Range("E45").Style = $origStyle
$origStyle = $ws.Range("E46").Style
$ws.Range("E46").Value = "'0.18%"
$ws.Range("E46").Style = $origStyle
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").Value = "'0.1091"
$ws.Range("D47").Style = $origStyle
$origStyle = $ws.Range("E47").Style
$ws.Range("E47").Value = "'-5.28%"
$ws.Range("E47").Style = $origStyle
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").Value = "'0.002143"
$ws.Range("D48").Style = $origStyle
$origStyle = $ws.Range("E48").Style
$ws.Range("E48").Value = "'-12.15%"
$ws.Range("E48").Style = $origStyle
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").Value = "'0.00002104"
$ws.Range("D49").Style = $origStyle
$origStyle = $ws.Range("E49").Style
$ws.Range("E49").Value = "'0.18%"
$ws.Range("E49").Style = $origStyle
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").Value = "'0.0002004"
$ws.Range("D50").Style = $origStyle
$origStyle = $ws.Range("E50").Style
$ws.Range("E50").Value = "'0.18%"
$ws.Range("E50").Style = $origStyle
